$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.654.69"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.301.62"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.03"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.32"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.80"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.34"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.27"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "2.651.48"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "2.301.31"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "42.623.77"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.39"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.61"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.40"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.71"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.86"
$ws.Range("E26").Value = "  +17.98%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.91"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.88"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.25"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.14"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0865"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.81"
$ws.Range("E34").Value = "  +6.34%  "
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.62"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.90"
$ws.Range("E41").Value = "  +10.34%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.00"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.26"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "1.697.76"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.57"
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.05"
$ws.Range("E49").Value = "  -4.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.83"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  -2.62%  "
